$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Internal Assignment" column (O) is dropped from the export: clear the
# header cell and the per-property "FALSE" values in both property-type
# tables (rows 12-15 and rows 20-23), keeping the existing cell styles
# (this also makes the now-unused "Internal Assignment" shared string
# disappear from the workbook on save).
$ws.Range("O12:O15").ClearContents() | Out-Null
$ws.Range("O20:O23").ClearContents() | Out-Null

# Row 17 (the second "SAMPLE_TYPE" section header) got a little shorter.
$ws.Rows.Item(17).RowHeight = 17.35

# Reflect the resulting selection/active cell.
$ws.Range("O12:O23").Select() | Out-Null
